$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Boundary_condition") before the old Scheme_order
# column, shifting Scheme_order, Mesh_cell_type, Test_color and
# Computational_time one column to the right.
$ws.Columns.Item(6).Insert()

# New header for the inserted column - copy the header formatting
# (bold font, border, centered alignment) from the neighboring header cell
$ws.Range("F1").Value = "Boundary_condition"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill in the Boundary_condition values for each data row
$ws.Range("F2").Value = "Dirichlet"
$ws.Range("F3").Value = "Dirichlet"
$ws.Range("F4").Value = "Dirichlet"
$ws.Range("F5").Value = "Dirichlet"
$ws.Range("F6").Value = "Dirichlet"
$ws.Range("F7").Value = "Neumann"
$ws.Range("F8").Value = "Dirichlet"
$ws.Range("F9").Value = "Neumann"
$ws.Range("F10").Value = "Dirichlet"
$ws.Range("F11").Value = "Dirichlet"
$ws.Range("F12").Value = "Dirichlet"
$ws.Range("F13").Value = "Dirichlet"
$ws.Range("F14").Value = "Dirichlet"

# Updated Computational_time values (now in column J after the shift)
$ws.Range("J2").Value = 109.4771571159363
$ws.Range("J3").Value = 6.233874082565308
$ws.Range("J4").Value = 208.5218908786774
$ws.Range("J5").Value = 12.18183302879333
$ws.Range("J6").Value = 10.15054607391357
$ws.Range("J7").Value = 9.985594987869263
$ws.Range("J8").Value = 18.92783284187317
$ws.Range("J9").Value = 18.37629008293152
$ws.Range("J10").Value = 4.93955397605896
$ws.Range("J11").Value = 2.616051912307739
$ws.Range("J12").Value = 5.918725967407227
$ws.Range("J13").Value = 64.21971893310547
$ws.Range("J14").Value = 3.755897998809814

Write-Output "done"
